# Updated symbol list on Thu Jan 26 19:32:15 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the coin rows on the active sheet. Values are written as plain text (not
# numbers) - NumberFormat is forced to "@" (Text) on each target cell first
# so Excel does not reinterpret the numeric-looking strings/percentages as
# numbers (which would change their stored precision/representation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "304.85" },
    @{ Cell = "E2";  Value = "1.07%" },
    @{ Cell = "D3";  Value = "35.90" },
    @{ Cell = "E3";  Value = "1.15%" },
    @{ Cell = "D4";  Value = "5.035" },
    @{ Cell = "E4";  Value = "-0.76%" },
    @{ Cell = "D5";  Value = "0.08078" },
    @{ Cell = "E5";  Value = "0.97%" },
    @{ Cell = "D6";  Value = "1.915" },
    @{ Cell = "E6";  Value = "-0.78%" },
    @{ Cell = "D7";  Value = "4.137" },
    @{ Cell = "E7";  Value = "2.36%" },
    @{ Cell = "D8";  Value = "7.844" },
    @{ Cell = "E8";  Value = "0.97%" },
    @{ Cell = "D9";  Value = "0.9314" },
    @{ Cell = "E9";  Value = "0.38%" },
    @{ Cell = "D10"; Value = "0.1257" },
    @{ Cell = "E10"; Value = "-19.15%" },
    @{ Cell = "D11"; Value = "0.1908" },
    @{ Cell = "E11"; Value = "0.57%" },
    @{ Cell = "D12"; Value = "0.09203" },
    @{ Cell = "E12"; Value = "2.24%" },
    @{ Cell = "D13"; Value = "0.03501" },
    @{ Cell = "E13"; Value = "2.03%" },
    @{ Cell = "D14"; Value = "0.09935" },
    @{ Cell = "E14"; Value = "0.45%" },
    @{ Cell = "D15"; Value = "0.001422" },
    @{ Cell = "E15"; Value = "1.48%" },
    @{ Cell = "D16"; Value = "0.006651" },
    @{ Cell = "E16"; Value = "15.66%" },
    @{ Cell = "E17"; Value = "2.39%" },
    @{ Cell = "D18"; Value = "3.245" },
    @{ Cell = "E18"; Value = "9.46%" },
    @{ Cell = "D19"; Value = "0.3443" },
    @{ Cell = "E19"; Value = "-0.06%" },
    @{ Cell = "D20"; Value = "5.185" },
    @{ Cell = "E20"; Value = "3.11%" },
    @{ Cell = "D21"; Value = "0.1305" },
    @{ Cell = "E21"; Value = "0.16%" },
    @{ Cell = "E22"; Value = "5.71%" },
    @{ Cell = "D24"; Value = "0.001235" },
    @{ Cell = "D25"; Value = "0.004721" },
    @{ Cell = "E25"; Value = "-1.05%" },
    @{ Cell = "D26"; Value = "0.0001302" },
    @{ Cell = "E26"; Value = "5.85%" },
    @{ Cell = "E27"; Value = "3.63%" },
    @{ Cell = "D39"; Value = "0.01964" },
    @{ Cell = "E39"; Value = "6.39%" },
    @{ Cell = "D40"; Value = "0.05162" },
    @{ Cell = "E40"; Value = "8.36%" },
    @{ Cell = "D41"; Value = "0.007597" },
    @{ Cell = "E41"; Value = "4.09%" },
    @{ Cell = "D42"; Value = "0.01018" },
    @{ Cell = "E42"; Value = "-3.96%" },
    @{ Cell = "D43"; Value = "0.1371" },
    @{ Cell = "E43"; Value = "3.10%" },
    @{ Cell = "E44"; Value = "-0.33%" },
    @{ Cell = "D45"; Value = "0.01068" },
    @{ Cell = "E45"; Value = "10.04%" },
    @{ Cell = "D46"; Value = "0.00006389" },
    @{ Cell = "E46"; Value = "2.53%" },
    @{ Cell = "E47"; Value = "0.15%" },
    @{ Cell = "D48"; Value = "63.57" },
    @{ Cell = "E48"; Value = "-1.70%" },
    @{ Cell = "E49"; Value = "-3.46%" },
    @{ Cell = "D50"; Value = "0.00002103" },
    @{ Cell = "E50"; Value = "0.15%" },
    @{ Cell = "D51"; Value = "0.0002003" },
    @{ Cell = "E51"; Value = "0.15%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
